# Rebuilt Jacis pathfinder and generated red11
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Column B width ---
$ws.Columns("B").ColumnWidth = 34.8

# --- Row 5: new "Red11 / Feed Switch from Front Start 1 aud" data, highlighted yellow ---
$ws.Range("A5:O5").Interior.Color = 65535
$ws.Range("B5").Value = "Feed Switch from Front Start 1 aud"
$ws.Range("C5").Value = 4
$ws.Range("D5").Value = 16
$ws.Range("E5").Value = 280.25
$ws.Range("F5").Value = 0
$ws.Range("G5").Value = 80
$ws.Range("H5").Value = 290.25
$ws.Range("I5").Value = 75
$ws.Range("J5").Value = 120
$ws.Range("K5").Value = 290.25
$ws.Range("L5").Value = -30
$ws.Range("M5").Value = 140
$ws.Range("N5").Value = 255
$ws.Range("O5").Value = -110

# --- Row 6: now holds what used to be row 5's data ("Feed Switch from Side start 1 aud") ---
$ws.Range("B6").Value = "Feed Switch from Side start 1 aud"
$ws.Range("C6").Value = 4
$ws.Range("G6").Value = 80
$ws.Range("H6").Value = 234
$ws.Range("I6").Value = 25
$ws.Range("J6").Value = 118
$ws.Range("K6").Value = 234
$ws.Range("L6").Value = 0

# --- Rows 18 and 19: highlight yellow (no value changes) ---
$ws.Range("A18:L19").Interior.Color = 65535

# --- Record sort state on A2:R32 sorted by column A ---
$sortObj = $ws.Sort
$sortObj.SortFields.Clear()
$sortObj.SortFields.Add($ws.Range("A2:A32"))
$sortObj.SetRange($ws.Range("A2:R32"))
$sortObj.Header = 0
$sortObj.Apply()

# --- View: scroll/selection changes ---
$ws.Range("H18").Select()
$excel.ActiveWindow.ScrollRow = 4
